$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.477.48"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.422.93"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.53"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.64"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.423.57"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.94"
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.416"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "4.020.85"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.134"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.60"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "66.529.66"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000173"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "3.424.23"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.93"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.84"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.70"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.92"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000127"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.537"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.84"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.99"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.54"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.06"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.98"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.78"
$ws.Range("E39").Value = "  -5.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "2.724.96"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.39"
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0692"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.46"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.14"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "333.75"
$ws.Range("E48").Value = "  +7.08%  "
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.88"
$ws.Range("E51").Value = "  +3.15%  "
